$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# NOTE on index <-> physical-part mapping in this document:
#   Headers.Item(1) (wdHeaderFooterPrimary) -> word/header2.xml
#   Headers.Item(2) (wdHeaderFooterFirstPage) -> word/header1.xml
#   Footers.Item(1) (wdHeaderFooterPrimary) -> word/footer2.xml
#   Footers.Item(2) (wdHeaderFooterFirstPage) -> word/footer1.xml
# (confirmed empirically against the package's word/_rels/document.xml.rels,
#  which wires rId10/"default" to header2.xml and rId11/"first" to header1.xml,
#  and similarly for the footers)

# Pearson logo pictures (footer1.xml & footer2.xml): image1.png -> image2.png
$sec.Footers.Item(1).Range.InlineShapes.Item(1).Name = "image2.png"
$sec.Footers.Item(2).Range.InlineShapes.Item(1).Name = "image2.png"

# BTec logo pictures (header1.xml & header2.xml): image2.jpg -> image1.jpg
$sec.Headers.Item(1).Range.InlineShapes.Item(1).Name = "image1.jpg"
$sec.Headers.Item(2).Range.InlineShapes.Item(1).Name = "image1.jpg"
